$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.361.87'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.880.00'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08052'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3160'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.46%  '
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08342'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.893.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.263'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7181'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.379'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008662'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.375.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.154.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.825'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1571'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.095'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.507'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.437'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('E32').Value = '  -6.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05413'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.945'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7730'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.189'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01890'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.270.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.750'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.518'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.61%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '113.34'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.36%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9126'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '74.69'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('E46').Value = '  +5.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.036.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.814'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5222'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('E50').Value = '  +1.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4384'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.39%  '
